$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to text so numeric-looking values
# (e.g. "1.002", "30.108.24") are stored as strings, not coerced to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value2 = '30.108.24'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value2 = '1.920.56'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('D4').Value2 = '1.002'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value2 = '319.43'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value2 = '1.001'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value2 = '0.5077'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').Value2 = '0.4049'
$ws.Range('E8').Value = '  +2.42%  '
$ws.Range('D9').Value2 = '0.08329'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').Value2 = '1.116'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').Value2 = '42.09'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value2 = '24.18'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value2 = '6.425'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').Value2 = '1.918.01'
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('D15').Value2 = '7.251'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value2 = '1.002'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value2 = '92.58'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value2 = '0.00001095'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').Value2 = '0.06499'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('D20').Value2 = '18.47'
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('D21').Value2 = '1.001'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value2 = '5.955'
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').Value2 = '30.125.95'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').Value2 = '11.35'
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').Value2 = '2.196'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').Value2 = '2.136.88'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').Value2 = '21.86'
$ws.Range('E27').Value = '  +3.26%  '
$ws.Range('D28').Value2 = '162.60'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').Value2 = '2.263'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('D30').Value2 = '129.13'
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').Value2 = '1.134'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('D32').Value2 = '0.1047'
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').Value2 = '5.950'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').Value2 = '3.793'
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('D35').Value2 = '0.02448'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').Value2 = '5.317'
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('D37').Value2 = '0.06452'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').Value2 = '1.234'
$ws.Range('E38').Value = '  +4.44%  '
$ws.Range('D39').Value2 = '0.2149'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value2 = '0.6469'
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').Value2 = '8.626'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value2 = '11.49'
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value2 = '1.214'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('D44').Value2 = '13.28'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('D45').Value2 = '0.6051'
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('D46').Value2 = '2.171'
$ws.Range('E46').Value = '  +7.14%  '
$ws.Range('D47').Value2 = '3.625'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').Value2 = '122.33'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value2 = '1.208'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value2 = '1.138'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').Value2 = '78.07'
$ws.Range('E51').Value = '  +0.96%  '

# Restore the original (default/general) cell style on the Price column
# now that the text values are committed, so no extra formatting persists.
$priceRange.Style = "Normal"

